$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1089.6
$ws.Range("I107").Value = 942.8570999999999
$ws.Range("K107").Value = 942.8570999999999
$ws.Range("M107").Value = 977.1429000000001
$ws.Range("H118").Value = 1169.1538
$ws.Range("I118").Value = 1035
$ws.Range("J118").Value = 1616.3334
$ws.Range("K118").Value = 3105
$ws.Range("L118").Value = 4849.0002
$ws.Range("M118").Value = -1448
$ws.Range("N118").Value = -8163.0002
$ws.Range("H131").Value = 1863.96
$ws.Range("I131").Value = 1351.6923
$ws.Range("J131").Value = 2418.9167
$ws.Range("K131").Value = 4055.0769
$ws.Range("L131").Value = 7256.750100000001
$ws.Range("M131").Value = 984.9231
$ws.Range("N131").Value = -17336.7501
$ws.Range("H132").Value = 19814.5
$ws.Range("I132").Value = 2499.739
$ws.Range("J132").Value = 218934.25
$ws.Range("K132").Value = 7499.217000000001
$ws.Range("L132").Value = 656802.75
$ws.Range("M132").Value = -4969.217000000001
$ws.Range("N132").Value = -661862.75
$ws.Range("H135").Value = 62504100
$ws.Range("I135").Value = 2799.1667
$ws.Range("K135").Value = 25192.5003
$ws.Range("M135").Value = -22657.5003
$ws.Range("H141").Value = 3999.5454
$ws.Range("I141").Value = 2114.2856
$ws.Range("K141").Value = 6342.8568
$ws.Range("M141").Value = -1162.8568
$ws.Range("H62").Value = 5169.375
$ws.Range("I62").Value = 5446.923
$ws.Range("K62").Value = 5446.923
$ws.Range("M62").Value = -4822.923
$ws.Range("H65").Value = 5169.375
$ws.Range("I65").Value = 5446.923
$ws.Range("K65").Value = 27234.615
$ws.Range("M65").Value = -24114.615

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1421.0454
$ws.Range("I110").Value = 1417.2858
$ws.Range("J110").Value = 1500
$ws.Range("K110").Value = 1417.2858
$ws.Range("L110").Value = 1500
$ws.Range("M110").Value = 627.7141999999999
$ws.Range("N110").Value = -5590
$ws.Range("H32").Value = 15494.5
$ws.Range("I32").Value = 14752.667
$ws.Range("J32").Value = 18991.715
$ws.Range("K32").Value = 14752.667
$ws.Range("L32").Value = 18991.715
$ws.Range("M32").Value = -14465.667
$ws.Range("N32").Value = -19565.715
$ws.Range("H63").Value = 4034.5
$ws.Range("I63").Value = 2593
$ws.Range("J63").Value = 5476
$ws.Range("K63").Value = 2593
$ws.Range("L63").Value = 5476
$ws.Range("M63").Value = -1907
$ws.Range("N63").Value = -6848
$ws.Range("H66").Value = 4034.5
$ws.Range("I66").Value = 2593
$ws.Range("J66").Value = 5476
$ws.Range("K66").Value = 12965
$ws.Range("L66").Value = 27380
$ws.Range("M66").Value = -9533
$ws.Range("N66").Value = -34244

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 73480
$ws.Range("J132").Value = 73480
$ws.Range("L132").Value = 73480
$ws.Range("N132").Value = -83600
$ws.Range("H134").Value = 2500.6438
$ws.Range("I134").Value = 1472.1111
$ws.Range("K134").Value = 4416.3333
$ws.Range("M134").Value = -1881.3333
$ws.Range("H35").Value = 35465
$ws.Range("J35").Value = 35465
$ws.Range("L35").Value = 35465
$ws.Range("N35").Value = -36085
$ws.Range("H82").Value = 15536.786
$ws.Range("I82").Value = 4378.75
$ws.Range("K82").Value = 4378.75
$ws.Range("M82").Value = -3995.75
$ws.Range("H85").Value = 15536.786
$ws.Range("I85").Value = 4378.75
$ws.Range("K85").Value = 4378.75
$ws.Range("M85").Value = -3052.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1412.3334
$ws.Range("J113").Value = 1200
$ws.Range("L113").Value = 1200
$ws.Range("N113").Value = -5540
$ws.Range("H136").Value = 2282.2964
$ws.Range("I136").Value = 1623.4667
$ws.Range("J136").Value = 3105.8333
$ws.Range("K136").Value = 4870.4001
$ws.Range("L136").Value = 9317.499899999999
$ws.Range("M136").Value = -2320.4001
$ws.Range("N136").Value = -14417.4999
$ws.Range("H16").Value = 1412.3334
$ws.Range("J16").Value = 1200
$ws.Range("L16").Value = 1200
$ws.Range("N16").Value = -1774
$ws.Range("H31").Value = 7582305
$ws.Range("I31").Value = 4887.875
$ws.Range("J31").Value = 9266175
$ws.Range("K31").Value = 4887.875
$ws.Range("L31").Value = 9266175
$ws.Range("M31").Value = -4592.875
$ws.Range("N31").Value = -9266765
$ws.Range("H34").Value = 7582305
$ws.Range("I34").Value = 4887.875
$ws.Range("J34").Value = 9266175
$ws.Range("K34").Value = 4887.875
$ws.Range("L34").Value = 9266175
$ws.Range("M34").Value = -4685.875
$ws.Range("N34").Value = -9266579
$ws.Range("H58").Value = 2282.2964
$ws.Range("I58").Value = 1623.4667
$ws.Range("J58").Value = 3105.8333
$ws.Range("K58").Value = 1623.4667
$ws.Range("L58").Value = 3105.8333
$ws.Range("M58").Value = -1420.4667
$ws.Range("N58").Value = -3511.8333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 13198.333
$ws.Range("J122").Value = 19584.166
$ws.Range("L122").Value = 176257.494
$ws.Range("N122").Value = -181157.494
$ws.Range("H129").Value = 121431
$ws.Range("I129").Value = 429437.16
$ws.Range("J129").Value = 1650.8334
$ws.Range("K129").Value = 1288311.48
$ws.Range("L129").Value = 4952.5002
$ws.Range("M129").Value = -1283311.48
$ws.Range("N129").Value = -14952.5002
$ws.Range("H135").Value = 5236.0454
$ws.Range("I135").Value = 5862.5264
$ws.Range("K135").Value = 52762.7376
$ws.Range("M135").Value = -50227.7376
$ws.Range("H5").Value = 5236.0454
$ws.Range("I5").Value = 5862.5264
$ws.Range("K5").Value = 17587.5792
$ws.Range("M5").Value = -17475.5792

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2706.9333
$ws.Range("I107").Value = 395.14285
$ws.Range("J107").Value = 4729.75
$ws.Range("K107").Value = 395.14285
$ws.Range("L107").Value = 4729.75
$ws.Range("M107").Value = 1524.85715
$ws.Range("N107").Value = -8569.75
$ws.Range("H113").Value = 1277.3334
$ws.Range("I113").Value = 1407.2727
$ws.Range("J113").Value = 920
$ws.Range("K113").Value = 1407.2727
$ws.Range("L113").Value = 920
$ws.Range("M113").Value = 762.7273
$ws.Range("N113").Value = -5260
$ws.Range("H122").Value = 1976
$ws.Range("I122").Value = 1781.7778
$ws.Range("J122").Value = 2850
$ws.Range("K122").Value = 5345.3334
$ws.Range("L122").Value = 8550
$ws.Range("M122").Value = -2895.3334
$ws.Range("N122").Value = -13450
$ws.Range("H132").Value = 2716.7954
$ws.Range("I132").Value = 1841.9546
$ws.Range("J132").Value = 3591.6365
$ws.Range("K132").Value = 5525.8638
$ws.Range("L132").Value = 10774.9095
$ws.Range("M132").Value = -2995.8638
$ws.Range("N132").Value = -15834.9095
$ws.Range("H38").Value = 16972.5
$ws.Range("J38").Value = 16972.5
$ws.Range("L38").Value = 16972.5
$ws.Range("N38").Value = -17898.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 80010
$ws.Range("J5").Value = 80010
$ws.Range("L5").Value = 80010
$ws.Range("N5").Value = -80236

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2220.2964
$ws.Range("I132").Value = 1566
$ws.Range("K132").Value = 4698
$ws.Range("M132").Value = -2168
$ws.Range("N39").ClearContents()
$ws.Range("H39").Value = 9999
$ws.Range("I39").Value = 9999
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 9999
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -9586
$ws.Range("H42").Value = 21742.5
$ws.Range("J42").Value = 21742.5
$ws.Range("L42").Value = 21742.5
$ws.Range("N42").Value = -22498.5
$ws.Range("H43").Value = 24489.25
$ws.Range("I43").Value = 11027
$ws.Range("J43").Value = 28976.666
$ws.Range("K43").Value = 11027
$ws.Range("L43").Value = 28976.666
$ws.Range("M43").Value = -10878
$ws.Range("N43").Value = -29274.666
